$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H15").Value = 1672.6552
$ws.Range("I15").Value = 1672.6552
$ws.Range("K15").Value = 5017.9656
$ws.Range("M15").Value = -4848.9656

$ws.Range("H43").Value = 1316.3334
$ws.Range("J43").Value = 900
$ws.Range("L43").Value = 900
$ws.Range("N43").Value = -1038

$ws.Range("H107").Value = 862.04346
$ws.Range("I107").Value = 892.2273
$ws.Range("K107").Value = 892.2273
$ws.Range("M107").Value = 1027.7727

$ws.Range("H132").Value = 1269.5834
$ws.Range("I132").Value = 1269.5834
$ws.Range("K132").Value = 3808.7502
$ws.Range("M132").Value = -1278.7502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8947.77
$ws.Range("I32").Value = 8592.056
$ws.Range("K32").Value = 8592.056
$ws.Range("M32").Value = -8305.056

$ws.Range("H97").Value = 3489.9
$ws.Range("I97").Value = 1199.8572
$ws.Range("K97").Value = 1199.8572
$ws.Range("M97").Value = -703.8571999999999

$ws.Range("H122").Value = 1495.1428
$ws.Range("I122").Value = 1225.5385
$ws.Range("K122").Value = 3676.6155
$ws.Range("M122").Value = -1226.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5154.1816
$ws.Range("I20").Value = 1862.125
$ws.Range("J20").Value = 13933
$ws.Range("K20").Value = 1862.125
$ws.Range("L20").Value = 13933
$ws.Range("M20").Value = -1615.125
$ws.Range("N20").Value = -14427

$ws.Range("H22").Value = 105.57143
$ws.Range("I22").Value = 77.8
$ws.Range("J22").Value = 175
$ws.Range("K22").Value = 77.8
$ws.Range("L22").Value = 175
$ws.Range("M22").Value = 95.2
$ws.Range("N22").Value = -521

$ws.Range("H107").Value = 1614.3334
$ws.Range("J107").Value = 1488.5
$ws.Range("L107").Value = 1488.5
$ws.Range("N107").Value = -5328.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10000375
$ws.Range("I6").Value = 20000000
$ws.Range("J6").Value = 750
$ws.Range("K6").Value = 20000000
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -19999887
$ws.Range("N6").Value = -976

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H25").Value = 40000
$ws.Range("J25").Value = 40000
$ws.Range("L25").Value = 40000
$ws.Range("N25").Value = -40348

$ws.Range("H31").Value = 1849.6428
$ws.Range("I31").Value = 1045
$ws.Range("J31").Value = 4800
$ws.Range("K31").Value = 1045
$ws.Range("L31").Value = 4800
$ws.Range("M31").Value = -750
$ws.Range("N31").Value = -5390

$ws.Range("H34").Value = 1849.6428
$ws.Range("I34").Value = 1045
$ws.Range("J34").Value = 4800
$ws.Range("K34").Value = 1045
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -843
$ws.Range("N34").Value = -5204

$ws.Range("H41").Value = 35000
$ws.Range("J41").Value = 35000
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35856

$ws.Range("H50").Value = 39957.5
$ws.Range("J50").Value = 39957.5
$ws.Range("L50").Value = 39957.5
$ws.Range("N50").Value = -41207.5

$ws.Range("H60").Value = 47724.5
$ws.Range("J60").Value = 47724.5
$ws.Range("L60").Value = 47724.5
$ws.Range("N60").Value = -48746.5

$ws.Range("H86").Value = 5056.5713
$ws.Range("I86").Value = 3996.2
$ws.Range("J86").Value = 6020.5454
$ws.Range("K86").Value = 3996.2
$ws.Range("L86").Value = 6020.5454
$ws.Range("M86").Value = -2873.2
$ws.Range("N86").Value = -8266.545399999999

$ws.Range("H89").Value = 5056.5713
$ws.Range("I89").Value = 3996.2
$ws.Range("J89").Value = 6020.5454
$ws.Range("K89").Value = 19981
$ws.Range("L89").Value = 30102.727
$ws.Range("M89").Value = -14365
$ws.Range("N89").Value = -41334.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("M16").Value = -427

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H55").Value = 13818.2
$ws.Range("J55").Value = 13818.2
$ws.Range("L55").Value = 41454.60000000001
$ws.Range("N55").Value = -41808.60000000001

$ws.Range("H70").Value = 2012
$ws.Range("I70").Value = 2012
$ws.Range("K70").Value = 6036
$ws.Range("M70").Value = -5721

$ws.Range("H73").Value = 2012
$ws.Range("I73").Value = 2012
$ws.Range("K73").Value = 6036
$ws.Range("M73").Value = -4944

$ws.Range("H107").Value = 442.5
$ws.Range("J107").Value = 685
$ws.Range("L107").Value = 2055
$ws.Range("N107").Value = -5895

$ws.Range("H121").Value = 2469.9
$ws.Range("I121").Value = 1499.5
$ws.Range("J121").Value = 2712.5
$ws.Range("K121").Value = 4498.5
$ws.Range("L121").Value = 8137.5
$ws.Range("M121").Value = -3188.5
$ws.Range("N121").Value = -10757.5

$ws.Range("H136").Value = 2363.6667
$ws.Range("I136").Value = 1977.5
$ws.Range("K136").Value = 5932.5
$ws.Range("M136").Value = -832.5

$ws.Range("H138").Value = 2450
$ws.Range("J138").Value = 5250
$ws.Range("L138").Value = 15750
$ws.Range("N138").Value = -26030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1573.7142
$ws.Range("I97").Value = 1671
$ws.Range("K97").Value = 1671
$ws.Range("M97").Value = -1175

$ws.Range("H122").Value = 2417
$ws.Range("I122").Value = 2827.1428
$ws.Range("J122").Value = 1842.8
$ws.Range("K122").Value = 8481.428400000001
$ws.Range("L122").Value = 5528.4
$ws.Range("M122").Value = -6031.428400000001
$ws.Range("N122").Value = -10428.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3311.8572
$ws.Range("I7").Value = 3311.8572
$ws.Range("K7").Value = 3311.8572
$ws.Range("M7").Value = -3199.8572

$ws.Range("H31").Value = 9338
$ws.Range("I31").Value = 5507
$ws.Range("J31").Value = 17000
$ws.Range("K31").Value = 5507
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = -5259
$ws.Range("N31").Value = -17496

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H55").Value = 1305.25
$ws.Range("J55").Value = 1324.5
$ws.Range("L55").Value = 1324.5
$ws.Range("N55").Value = -1670.5

$ws.Range("H122").Value = 3453.2
$ws.Range("I122").Value = 3453.2
$ws.Range("K122").Value = 10359.6
$ws.Range("M122").Value = -7909.599999999999

$ws.Range("H126").Value = 3311.8572
$ws.Range("I126").Value = 3311.8572
$ws.Range("K126").Value = 9935.571599999999
$ws.Range("M126").Value = -7465.571599999999

$ws.Range("H136").Value = 3377.8462
$ws.Range("I136").Value = 2300.1428
$ws.Range("K136").Value = 6900.428400000001
$ws.Range("M136").Value = -4350.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7821.4287
$ws.Range("I62").Value = 3966.6667
$ws.Range("J62").Value = 10712.5
$ws.Range("K62").Value = 3966.6667
$ws.Range("L62").Value = 10712.5
$ws.Range("M62").Value = -3342.6667
$ws.Range("N62").Value = -11960.5

$ws.Range("H65").Value = 7821.4287
$ws.Range("I65").Value = 3966.6667
$ws.Range("J65").Value = 10712.5
$ws.Range("K65").Value = 19833.3335
$ws.Range("L65").Value = 53562.5
$ws.Range("M65").Value = -16713.3335
$ws.Range("N65").Value = -59802.5

$ws.Range("H122").Value = 1075.1428
$ws.Range("I122").Value = 1075.1428
$ws.Range("K122").Value = 3225.4284
$ws.Range("M122").Value = -775.4284000000002
